# Weekly update: insert a new week's worth of "Repollo" price records
# (Vega Central Mapocho de Santiago) ahead of the existing historical
# rows, pushing the old rows 809-819 down to 813-823.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 blank rows at row 809 (old rows 809-819 shift to 813-823).
$ws.Range("A809:R812").EntireRow.Insert()

# Helper to fill one data row for this product block (positional args only
# -- named/typed parameters are not reliably supported by this runtime).
function Set-RepolloRow($Row, $Fecha, $Variedad, $Calidad, $Volumen, $PrecioMinimo, $PrecioMaximo, $PrecioPromedio, $Region) {
    $ws.Cells.Item($Row, 1).Value = 9
    $ws.Cells.Item($Row, 2).Value = "Vega Central Mapocho de Santiago"
    $ws.Cells.Item($Row, 3).Value = "Metropolitana"
    $ws.Cells.Item($Row, 4).Value = $Fecha
    $ws.Cells.Item($Row, 5).Value = 13
    $ws.Cells.Item($Row, 6).Value = 100112006
    $ws.Cells.Item($Row, 7).Value = "Repollo"
    $ws.Cells.Item($Row, 8).Value = $Variedad
    $ws.Cells.Item($Row, 9).Value = $Calidad
    $ws.Cells.Item($Row, 10).Value = $Volumen
    $ws.Cells.Item($Row, 11).Value = $PrecioMinimo
    $ws.Cells.Item($Row, 12).Value = $PrecioMaximo
    $ws.Cells.Item($Row, 13).Value = $PrecioPromedio
    $ws.Cells.Item($Row, 14).Value = "$/unidad"
    $ws.Cells.Item($Row, 15).Value = $Region
    $ws.Cells.Item($Row, 16).Value = $PrecioPromedio
    $ws.Cells.Item($Row, 17).Value = 1
    $ws.Cells.Item($Row, 18).Value = "Hortaliza"
}

Set-RepolloRow 809 45121 "Crespo record" "Primera" 2500 700 800 750 "Provincia de Quillota"
Set-RepolloRow 810 45121 "Crespo record" "Primera" 970 900 1000 948 "Región de O'Higgins"
Set-RepolloRow 811 45121 "Crespo record" "Segunda" 1690 600 600 600 "Provincia de Quillota"
Set-RepolloRow 812 45121 "Morada(o)" "Primera" 1600 900 1000 950 "Provincia de Quillota"
